$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Replacement paragraph XML for each of the 4 fields, in document order.
# Each field (a Word "m:..." instruction field, e.g. "{m: OrderedSet{...}->myTemplate()}")
# is converted from a field (fldChar begin/instrText/fldChar end) into plain literal
# text runs reading "{" + trimmed field code + "}", keeping the same run splits that the
# original instrText runs used.
$replacements = @(
    ('<w:p ' + $wns + ' w:rsidR="00735354" w:rsidRPr="00DC5685" w:rsidRDefault="00735354" w:rsidP="00F5495F">' +
        '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{m: </w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>OrderedSet</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{1, 2, 3, 4, 5}-&gt;</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">myTemplate()}</w:t></w:r>' +
     '</w:p>'),
    ('<w:p ' + $wns + ' w:rsidR="00735354" w:rsidRDefault="00735354" w:rsidP="00735354">' +
        '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{m:template myTemplate(a:</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>OrderedSet</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(Integer)</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">)}</w:t></w:r>' +
     '</w:p>'),
    ('<w:p ' + $wns + ' w:rsidR="00735354" w:rsidRDefault="00735354" w:rsidP="00735354">' +
        '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{m: a</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-&gt;sep(''['', ''|'', '']'')</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">}</w:t></w:r>' +
     '</w:p>'),
    ('<w:p ' + $wns + ' w:rsidR="007A2DC4" w:rsidRPr="00DC5685" w:rsidRDefault="00735354">' +
        '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{m:endtemplate}</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
     '</w:p>')
)

# Walk the paragraphs collecting those that contain a field (from first to last), and
# replace them one at a time. We go in reverse document order so that replacing a
# paragraph's range does not invalidate the character offsets of paragraphs that come
# after it (InsertXML can change the overall document length).
$fieldParas = New-Object System.Collections.ArrayList
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        [void]$fieldParas.Add($p)
    }
}

for ($idx = $fieldParas.Count - 1; $idx -ge 0; $idx--) {
    $p = $fieldParas[$idx]
    $rng = $p.Range
    [void]$rng.InsertXML($replacements[$idx])
}
